$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "UMNpw6rfS2pmonOusb6e"
$ws.Range("D3").Value = "ghfk30qKrgCNRBSKlCI6"
$ws.Range("E3").Value = "fBYZuVBGUyS04DphfKX5"
$ws.Range("F3").Value = "8ChZgJJQZjwBGQ8eTbt0"

$ws.Range("C9").Value = "7ecmDZLun7BADuXpFgD4"
$ws.Range("D9").Value = "juTbXWy1B7bfDcPBzxoU"
$ws.Range("E9").Value = "3tII00GFHRtOlrvGdOJB"
$ws.Range("F9").Value = "kZ90523POQNyAy9ozQxQ"

$ws.Range("F9").Select()
